$title42 = @'
"Association of
            <scp>HLA</scp>
            gene polymorphism with susceptibility, severity, and mortality of
            <scp>COVID</scp>
            ‐19: A systematic review"
'@
$authors48 = @'
[Paroma%Deb%xref no email%1,  Kaniz‐e‐%Zannat%xref no email%1,  Shiny%Talukder%xref no email%1,  Amirul Huda%Bhuiyan%xref no email%1,  Md. Shariful Alam%Jilani%xref no email%1,  K. M.%Saif‐Ur‐Rahman%xref no email%1]
'@
$idformat44 = @'
CROSSREF
'@
$dateacc45 = @'
2023-05-18
'@
$authors49 = @'
[Hussein N.%Ali%NULL%1,        Sherko S.%Niranji%sherko.subhan@garmian.edu.krd%1,        Sirwan M. A.%Al‐Jaf%NULL%2,        Sirwan M. A.%Al‐Jaf%NULL%0]
'@
$authors50 = @'
[Ana Valesca Fernandes Gilson%Silva%NULL%1,        Diego%Menezes%NULL%1,        Filipe Romero Rebello%Moreira%NULL%1,        Octávio Alcântara%Torres%NULL%1,        Paula Luize Camargos%Fonseca%NULL%1,        Rennan Garcias%Moreira%NULL%1,        Hugo José%Alves%NULL%1,        Vivian Ribeiro%Alves%NULL%1,        Tânia Maria de Resende%Amaral%NULL%1,        Adriano Neves%Coelho%NULL%1,        Júlia Maria%Saraiva Duarte%NULL%1,        Augusto Viana%da Rocha%NULL%1,        Luiz Gonzaga Paula%de Almeida%NULL%1,        João Locke Ferreira%de Araújo%NULL%1,        Hilton Soares%de Oliveira%NULL%1,        Nova Jersey Cláudio%de Oliveira%NULL%1,        Camila%Zolini%NULL%1,        Jôsy Hubner%de Sousa%NULL%1,        Elizângela Gonçalves%de Souza%NULL%1,        Rafael Marques%de Souza%NULL%1,        Luciana de Lima%Ferreira%NULL%1,        Alexandra%Lehmkuhl Gerber%NULL%1,        Ana Paula de Campos%Guimarães%NULL%1,        Paulo Henrique Silva%Maia%NULL%1,        Fernanda Martins%Marim%NULL%1,        Lucyene%Miguita%NULL%1,        Cristiane Campos%Monteiro%NULL%1,        Tuffi Saliba%Neto%NULL%1,        Fabrícia Soares Freire%Pugêdo%NULL%1,        Daniel Costa%Queiroz%NULL%1,        Damares Nigia Alborguetti Cuzzuol%Queiroz%NULL%1,        Luciana Cunha%Resende-Moreira%NULL%1,        Franciele Martins%Santos%NULL%1,        Erika Fernanda Carlos%Souza%NULL%1,        Carolina Moreira%Voloch%NULL%1,        Ana Tereza%Vasconcelos%NULL%1,        Renato Santana%de Aguiar%NULL%1,        Renan Pedra%de Souza%NULL%1]
'@
$doi2 = @'
10.1111/tan.14560
'@

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (10.1111/tan.14560) -- Date Accepted (H2) is handled first, before the
# other row-2 writes, so the multi-line title in C2 can't influence Excel's
# row-format recalculation while we coerce H2 back to a plain/default style.
# Force text storage so the date-looking string isn't auto-converted to a date
# serial number, then strip the temporary number format back off so the cell
# keeps its original (default) style.
$ws.Cells.Item(2, 8).NumberFormat = "@"
$ws.Cells.Item(2, 8).Value = $dateacc45
$ws.Cells.Item(2, 8).ClearFormats()

# Row 2 (10.1111/tan.14560) -- Title, Authors, ID, ID Format
$ws.Cells.Item(2, 3).Value = $title42
$ws.Cells.Item(2, 5).Value = $authors48
$ws.Cells.Item(2, 6).Value = $doi2
$ws.Cells.Item(2, 7).Value = $idformat44

# Row 3 (10.1002/jcla.24400) -- Authors whitespace-variant update
$ws.Cells.Item(3, 5).Value = $authors49

# Row 4 (10.3389/fmicb.2022.799713) -- Authors whitespace-variant update
$ws.Cells.Item(4, 5).Value = $authors50

# The new multi-line title in C2 makes Excel auto-grow row 2's height; put it
# back to the sheet's implicit auto height so no stray row-height override is
# written out.
$ws.Rows(2).AutoFit()
